$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Taakstructuur-referenties toevoegen/verplaatsen in de teamkalender

# "personal touch & bieding" verhuist van week 4 (rij 9) naar week 3 (rij 8)
$ws.Range("L9").Value = ""
$ws.Range("L8").Value = "personal touch & bieding"

# "tussentijds verslag" krijgt sectieverwijzing [5.1]
$ws.Range("H12").Value = "tussentijds verslag [5.1]"

# "eindpresentatie en verslag" krijgt sectieverwijzingen [5.2, 5.3]
$ws.Range("H19").Value = "eindpresentatie en verslag [5.2, 5.3]"

# "tests wagen [4]" wordt "tests wagen [4.2]"
$ws.Range("H16").Value = "tests wagen [4.2]"

# nieuwe taak "test zonder parcour [4.1]" toegevoegd in week 9
$ws.Range("H15").Value = "test zonder parcour [4.1]"

# selectie verplaatst naar H18
$ws.Range("H18").Select()
